$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force every written cell to remain Text (matches source inlineStr cells),
# since plain numeric-looking strings would otherwise be auto-converted to numbers.
function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '43.542.69'

Set-TextValue "D3" '2.416.48'
Set-TextValue "E3" '  +2.65%  '

Set-TextValue "E4" '  +0.10%  '

Set-TextValue "D5" '306.29'
Set-TextValue "E5" '  +1.17%  '

Set-TextValue "D6" '97.66'
Set-TextValue "E6" '  +2.08%  '

Set-TextValue "E7" '  +0.61%  '

Set-TextValue "E8" '  +0.06%  '

Set-TextValue "D9" '0.491'
Set-TextValue "E9" '  -1.60%  '

Set-TextValue "D10" '35.04'
Set-TextValue "E10" '  +2.48%  '

Set-TextValue "E11" '  +3.15%  '

Set-TextValue "D12" '0.0796'
Set-TextValue "E12" '  +0.88%  '

Set-TextValue "D13" '18.59'
Set-TextValue "E13" '  -0.67%  '

Set-TextValue "D14" '6.88'
Set-TextValue "E14" '  +2.01%  '

Set-TextValue "D15" '2.786.76'
Set-TextValue "E15" '  +2.36%  '

Set-TextValue "D16" '2.462.99'
Set-TextValue "E16" '  +6.13%  '

Set-TextValue "D17" '0.827'
Set-TextValue "E17" '  +3.69%  '

Set-TextValue "D18" '43.613.67'
Set-TextValue "E18" '  +1.02%  '

Set-TextValue "D19" '12.20'
Set-TextValue "E19" '  -0.16%  '

Set-TextValue "D20" '6.43'
Set-TextValue "E20" '  +2.69%  '

Set-TextValue "D21" '0.0₃0903'
Set-TextValue "E21" '  +1.27%  '

Set-TextValue "D22" '68.66'
Set-TextValue "E22" '  +0.73%  '

Set-TextValue "D23" '237.99'
Set-TextValue "E23" '  +0.89%  '

Set-TextValue "D24" '2.24'
Set-TextValue "E24" '  +0.74%  '

Set-TextValue "D25" '2.45'
Set-TextValue "E25" '  +0.78%  '

Set-TextValue "E26" '  +0.03%  '

Set-TextValue "D27" '24.99'
Set-TextValue "E27" '  +1.77%  '

Set-TextValue "D28" '2.23'
Set-TextValue "E28" '  -5.57%  '

Set-TextValue "D29" '9.40'
Set-TextValue "E29" '  +2.85%  '

Set-TextValue "D30" '32.49'
Set-TextValue "E30" '  +3.36%  '

Set-TextValue "D31" '5.13'
Set-TextValue "E31" '  +1.97%  '

Set-TextValue "D32" '18.39'
Set-TextValue "E32" '  +6.99%  '

Set-TextValue "D33" '0.114'
Set-TextValue "E33" '  +13.74%  '

Set-TextValue "D34" '0.999'
Set-TextValue "E34" '  -0.09%  '

Set-TextValue "D35" '0.0744'
Set-TextValue "E35" '  +2.54%  '

Set-TextValue "D36" '133.21'
Set-TextValue "E36" '  +20.21%  '

Set-TextValue "B37" 'LidoDAOToken'
Set-TextValue "C37" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D37" '2.97'
Set-TextValue "E37" '  +7.50%  '

Set-TextValue "B38" 'ARBITRUM'
Set-TextValue "C38" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D38" '1.88'
Set-TextValue "E38" '  +2.41%  '

Set-TextValue "D39" '4.40'
Set-TextValue "E39" '  -0.08%  '

Set-TextValue "E40" '  -1.41%  '

Set-TextValue "E41" '  -0.17%  '

Set-TextValue "D42" '21.48'
Set-TextValue "E42" '  -4.43%  '

Set-TextValue "D43" '1.951.35'
Set-TextValue "E43" '  +0.36%  '

Set-TextValue "E44" '  +1.11%  '

Set-TextValue "D45" '2.16'
Set-TextValue "E45" '  +2.00%  '

Set-TextValue "E46" '  +2.41%  '

Set-TextValue "D47" '9.28'
Set-TextValue "E47" '  -0.97%  '

Set-TextValue "D48" '2.636.54'
Set-TextValue "E48" '  +2.01%  '

Set-TextValue "E49" '  +3.08%  '

Set-TextValue "D50" '52.58'
Set-TextValue "E50" '  -0.74%  '

Set-TextValue "D51" '72.25'
Set-TextValue "E51" '  +0.14%  '
